$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema6d"
$ws.Cells.Item(2, 3).Value = "Trem2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 24.44575933333333
$ws.Cells.Item(2, 8).Value = 73.337278
$ws.Cells.Item(2, 9).Value = 0.4034052273345712
$ws.Cells.Item(2, 10).Value = 0.4034052273345712
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 8.781144333333334
$ws.Cells.Item(2, 14).Value = 26.343433
$ws.Cells.Item(2, 15).Value = 0.1283382517649557
$ws.Cells.Item(2, 16).Value = 0.1283382517649557
$ws.Cells.Item(2, 17).Value = 214.6617410439304
$ws.Cells.Item(2, 18).Value = 1931.955669395374
$ws.Cells.Item(2, 19).Value = 0.05177232162896338
$ws.Cells.Item(2, 20).Value = 0.05177232162896337

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema6d"
$ws.Cells.Item(3, 3).Value = "Trem2"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 24.44575933333333
$ws.Cells.Item(3, 8).Value = 73.337278
$ws.Cells.Item(3, 9).Value = 0.4034052273345712
$ws.Cells.Item(3, 10).Value = 0.4034052273345712
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 59.64073466666667
$ws.Cells.Item(3, 14).Value = 178.922204
$ws.Cells.Item(3, 15).Value = 0.8716617482350444
$ws.Cells.Item(3, 16).Value = 0.8716617482350443
$ws.Cells.Item(3, 17).Value = 1457.963046124524
$ws.Cells.Item(3, 18).Value = 13121.66741512071
$ws.Cells.Item(3, 19).Value = 0.3516329057056079
$ws.Cells.Item(3, 20).Value = 0.3516329057056078

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Sema6d"
$ws.Cells.Item(4, 3).Value = "Trem2"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 14.28901333333333
$ws.Cells.Item(4, 8).Value = 42.86704
$ws.Cells.Item(4, 9).Value = 0.2357980618855278
$ws.Cells.Item(4, 10).Value = 0.2357980618855278
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 8.781144333333334
$ws.Cells.Item(4, 14).Value = 26.343433
$ws.Cells.Item(4, 15).Value = 0.1283382517649557
$ws.Cells.Item(4, 16).Value = 0.1283382517649557
$ws.Cells.Item(4, 17).Value = 125.4738884609245
$ws.Cells.Item(4, 18).Value = 1129.26499614832
$ws.Cells.Item(4, 19).Value = 0.03026191103195347
$ws.Cells.Item(4, 20).Value = 0.03026191103195347

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Sema6d"
$ws.Cells.Item(5, 3).Value = "Trem2"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 14.28901333333333
$ws.Cells.Item(5, 8).Value = 42.86704
$ws.Cells.Item(5, 9).Value = 0.2357980618855278
$ws.Cells.Item(5, 10).Value = 0.2357980618855278
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 59.64073466666667
$ws.Cells.Item(5, 14).Value = 178.922204
$ws.Cells.Item(5, 15).Value = 0.8716617482350444
$ws.Cells.Item(5, 16).Value = 0.8716617482350443
$ws.Cells.Item(5, 17).Value = 852.2072528617957
$ws.Cells.Item(5, 18).Value = 7669.865275756161
$ws.Cells.Item(5, 19).Value = 0.2055361508535744
$ws.Cells.Item(5, 20).Value = 0.2055361508535744

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Sema6d"
$ws.Cells.Item(6, 3).Value = "Trem2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.480005
$ws.Cells.Item(6, 8).Value = 1.440015
$ws.Cells.Item(6, 9).Value = 0.007921068169999337
$ws.Cells.Item(6, 10).Value = 0.007921068169999336
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 8.781144333333334
$ws.Cells.Item(6, 14).Value = 26.343433
$ws.Cells.Item(6, 15).Value = 0.1283382517649557
$ws.Cells.Item(6, 16).Value = 0.1283382517649557
$ws.Cells.Item(6, 17).Value = 4.214993185721667
$ws.Cells.Item(6, 18).Value = 37.934938671495
$ws.Cells.Item(6, 19).Value = 0.001016576041048752
$ws.Cells.Item(6, 20).Value = 0.001016576041048751

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Sema6d"
$ws.Cells.Item(7, 3).Value = "Trem2"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.480005
$ws.Cells.Item(7, 8).Value = 1.440015
$ws.Cells.Item(7, 9).Value = 0.007921068169999337
$ws.Cells.Item(7, 10).Value = 0.007921068169999336
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 59.64073466666667
$ws.Cells.Item(7, 14).Value = 178.922204
$ws.Cells.Item(7, 15).Value = 0.8716617482350444
$ws.Cells.Item(7, 16).Value = 0.8716617482350443
$ws.Cells.Item(7, 17).Value = 28.62785084367334
$ws.Cells.Item(7, 18).Value = 257.65065759306
$ws.Cells.Item(7, 19).Value = 0.006904492128950586
$ws.Cells.Item(7, 20).Value = 0.006904492128950584

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Sema6d"
$ws.Cells.Item(8, 3).Value = "Trem2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 21.38374133333333
$ws.Cells.Item(8, 8).Value = 64.151224
$ws.Cells.Item(8, 9).Value = 0.3528756426099016
$ws.Cells.Item(8, 10).Value = 0.3528756426099016
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 8.781144333333334
$ws.Cells.Item(8, 14).Value = 26.343433
$ws.Cells.Item(8, 15).Value = 0.1283382517649557
$ws.Cells.Item(8, 16).Value = 0.1283382517649557
$ws.Cells.Item(8, 17).Value = 187.7737190346658
$ws.Cells.Item(8, 18).Value = 1689.963471311992
$ws.Cells.Item(8, 19).Value = 0.04528744306299008
$ws.Cells.Item(8, 20).Value = 0.04528744306299008

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Sema6d"
$ws.Cells.Item(9, 3).Value = "Trem2"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 21.38374133333333
$ws.Cells.Item(9, 8).Value = 64.151224
$ws.Cells.Item(9, 9).Value = 0.3528756426099016
$ws.Cells.Item(9, 10).Value = 0.3528756426099016
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 59.64073466666667
$ws.Cells.Item(9, 14).Value = 178.922204
$ws.Cells.Item(9, 15).Value = 0.8716617482350444
$ws.Cells.Item(9, 16).Value = 0.8716617482350443
$ws.Cells.Item(9, 17).Value = 1275.342043041966
$ws.Cells.Item(9, 18).Value = 11478.0783873777
$ws.Cells.Item(9, 19).Value = 0.3075881995469116
$ws.Cells.Item(9, 20).Value = 0.3075881995469115
